$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 555.9727745343941
$ws.Range("D2").Value = 136.4468106118986
$ws.Range("F2").Value = 452
$ws.Range("G2").Value = 518
$ws.Range("H2").Value = 626
$ws.Range("I2").Value = 1993
$ws.Range("C3").Value = 40.66421552496731
$ws.Range("D3").Value = 4.785040199674193
$ws.Range("F3").Value = 37.69
$ws.Range("G3").Value = 39.94
$ws.Range("H3").Value = 43.35
$ws.Range("C4").Value = 1.386957574319462
$ws.Range("D4").Value = 2.251972007438511
$ws.Range("F4").Value = 0.51
$ws.Range("G4").Value = 1.01
$ws.Range("H4").Value = 1.79
$ws.Range("C5").Value = 323.8454538002499
$ws.Range("D5").Value = 10.23904542492872
$ws.Range("F5").Value = 317.99
$ws.Range("G5").Value = 325.67
$ws.Range("H5").Value = 332.25
$ws.Range("I5").Value = 342.81
$ws.Range("C6").Value = 21.03396809739433
$ws.Range("D6").Value = 2.100131944092435
$ws.Range("E6").Value = 15.22
$ws.Range("F6").Value = 19.48
$ws.Range("G6").Value = 20.66
$ws.Range("H6").Value = 22.2
$ws.Range("C7").Value = -76.31486173511111
$ws.Range("D7").Value = 22.43087720017906
$ws.Range("C8").Value = 7.772729262699947
$ws.Range("D8").Value = 6.826958247165576
$ws.Range("C9").Value = 9.322383399976859
$ws.Range("D9").Value = 1.688459844100303
$ws.Range("C10").Value = 867.8303429706599
$ws.Range("D10").Value = 0.4610121929352152
$ws.Range("C11").Value = 0.5569302727017788
$ws.Range("D11").Value = 0.5906666965420388
$ws.Range("C12").Value = 22.69027919682621
$ws.Range("D12").Value = 12.27804715606286
$ws.Range("C13").Value = 0.6716804562221651
$ws.Range("D13").Value = 0.7483016417848986
$ws.Range("C14").Value = 1.826072163270429
$ws.Range("D14").Value = 1.665922424560786
$ws.Range("C15").Value = 93.71486173511092
$ws.Range("D15").Value = 22.43087720014354
$ws.Range("C16").Value = -85.63531770634948
$ws.Range("D16").Value = 20.23553674694372
$ws.Range("F16").Value = -101.0778545523916
$ws.Range("G16").Value = -85.69305820175224
$ws.Range("H16").Value = -67.59612087980607
$ws.Range("C17").Value = -77.86258844364953
$ws.Range("D17").Value = 24.79941101361801
$ws.Range("F17").Value = -92.22214159641585
$ws.Range("G17").Value = -75.43249407632486
$ws.Range("H17").Value = -56.66683163887967
